# Update profiles, oysters and clams datasets + website
# Applies data corrections (pH / salinity unit fixes) on sheets L1, L2, L3, L4
# plus the resulting selection / active-sheet UI-state changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# L1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("L1")
$ws.Activate()

$ws.Range("C6").Value  = 7.91
$ws.Range("D6").Value  = 52.914999999999999
$ws.Range("C7").Value  = 7.88
$ws.Range("D7").Value  = 52.914000000000001
$ws.Range("C8").Value  = 7.99
$ws.Range("D8").Value  = 52.915999999999997
$ws.Range("C9").Value  = 8.15
$ws.Range("D9").Value  = 52.915999999999997
$ws.Range("C10").Value = 8.33
$ws.Range("D10").Value = 52.914999999999999
$ws.Range("C11").Value = 8.41
$ws.Range("D11").Value = 52.87
$ws.Range("C12").Value = 8.69
$ws.Range("D12").Value = 52.853000000000002
$ws.Range("C13").Value = 8.82
$ws.Range("D13").Value = 52.792999999999999
$ws.Range("C14").Value = 8.75
$ws.Range("D14").Value = 52.462000000000003

# C18:C26 values reversed (same data, rows re-ordered)
$ws.Range("C18").Value = 0.85
$ws.Range("C19").Value = 6.2
$ws.Range("C20").Value = 6.83
$ws.Range("C21").Value = 7.8
$ws.Range("C23").Value = 8.1999999999999993
$ws.Range("C24").Value = 8.75
$ws.Range("C25").Value = 8.73
$ws.Range("C26").Value = 8.7200000000000006

$ws.Range("E25").Select()

# ---------------------------------------------------------------------
# L2
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("L2")
$ws.Activate()

$ws.Range("C6").Value  = 8.09
$ws.Range("D6").Value  = 52.877000000000002
$ws.Range("C7").Value  = 8.09
$ws.Range("D7").Value  = 52.9
$ws.Range("C8").Value  = 8.24
$ws.Range("D8").Value  = 52.914999999999999
$ws.Range("C9").Value  = 8.42
$ws.Range("D9").Value  = 52.914999999999999
$ws.Range("C10").Value = 8.58
$ws.Range("D10").Value = 52.924999999999997
$ws.Range("C11").Value = 8.6
$ws.Range("D11").Value = 52.945
$ws.Range("C12").Value = 8.6999999999999993
$ws.Range("D12").Value = 52.942
$ws.Range("C13").Value = 8.84
$ws.Range("D13").Value = 52.88
$ws.Range("C14").Value = 8.9499999999999993
$ws.Range("D14").Value = 52.774999999999999

$ws.Range("E15").Select()

# ---------------------------------------------------------------------
# L3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("L3")
$ws.Activate()

$ws.Range("C6").Value  = 7.77
$ws.Range("D6").Value  = 52.95
$ws.Range("C7").Value  = 7.82
$ws.Range("D7").Value  = 52.966000000000001
$ws.Range("C8").Value  = 8.3000000000000007
$ws.Range("D8").Value  = 52.94
$ws.Range("C9").Value  = 8.36
$ws.Range("D9").Value  = 52.93
$ws.Range("C10").Value = 8.57
$ws.Range("D10").Value = 48.3

$ws.Range("D10").Select()

# ---------------------------------------------------------------------
# L4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("L4")
$ws.Activate()

$ws.Range("C6").Value = 9.1199999999999992
$ws.Range("D6").Value = 50.3
$ws.Range("C7").Value = 9.09
$ws.Range("D7").Value = 50

$ws.Range("F9").Select()
